# Auto-generated script applying the cryptos.xlsx price/volume update
# (GitHub Actions scheduled refresh), including the Polkadot/Chainlink row swap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.267.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.305.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.06%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.301.94"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.15%  "
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.24%  "
$ws.Range("E11").Value = "  -4.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.374"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.870.80"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.03%  "
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.307.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.07%  "
$ws.Range("E16").Value = "  -5.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.312.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "24.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.73%  "
$ws.Range("E21").Value = "  -10.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "354.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.76%  "
$ws.Range("E23").Value = "  -4.00%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.437.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.61%  "
$ws.Range("E27").Value = "  -5.33%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E31").Value = "  -5.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.37%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.332.68"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.13%  "
$ws.Range("E36").Value = "  -1.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "162.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("E40").Value = "  -3.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0752"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.30%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.742"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.31%  "
$ws.Range("E46").Value = "  -1.65%  "
$ws.Range("E47").Value = "  -4.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.89%  "
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.851"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.63%  "
